$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.615.26'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '1.792.36'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '231.78'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('D6').Value = '0.5878'
$ws.Range('E6').Value = '  -2.34%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.2764'
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.06744'
$ws.Range('E9').Value = '  -4.16%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '23.17'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '0.07527'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '1.797.03'
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('D13').Value = '4.783'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '0.6127'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '2.035.36'
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('E16').Value = '  -4.60%  '
$ws.Range('D17').Value = '0.000008885'
$ws.Range('E17').Value = '  -8.70%  '
$ws.Range('D18').Value = '28.595.70'
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('E19').Value = '  -6.98%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '209.03'
$ws.Range('E21').Value = '  -6.62%  '
$ws.Range('E22').Value = '  -1.96%  '
$ws.Range('D23').Value = '6.829'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '152.81'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '8.135'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('E27').Value = '  -2.99%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').Value = '1.416'
$ws.Range('E29').Value = '  -3.59%  '
$ws.Range('D30').Value = '0.06233'
$ws.Range('E30').Value = '  -5.52%  '
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').Value = '3.804'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = '3.781'
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').Value = '1.046'
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('D36').Value = '0.6385'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('D39').Value = '6.385'
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('D40').Value = '0.01694'
$ws.Range('E40').Value = '  -3.63%  '
$ws.Range('D41').Value = '1.141.01'
$ws.Range('E41').Value = '  -5.83%  '
$ws.Range('D42').Value = '0.8774'
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '1.943.90'
$ws.Range('E45').Value = '  -2.42%  '
$ws.Range('D46').Value = '59.87'
$ws.Range('E46').Value = '  -4.25%  '
$ws.Range('D47').Value = '0.00000000111'
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('D48').Value = '1.586'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '0.05466'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').Value = '8.315'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').Value = '0.4486'
$ws.Range('E51').Value = '  -1.44%  '
